# Update "Lương" sheet: remove the "Phụ cấp tại ..." rows for CẦN THƠ and
# LONG XUYÊN (these allowances no longer apply at those locations), shifting
# all subsequent rows up. SÓC TRĂNG keeps its "Phụ cấp" row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Lương")

# Delete "Phụ cấp tại LONG XUYÊN" first (row 14), then "Phụ cấp tại CẦN THƠ"
# (row 3), deleting the later row first keeps the earlier row index valid.
$ws.Rows.Item(14).Delete()
$ws.Rows.Item(3).Delete()
